$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for this market; insert it as a
# new row 7 (pushing the existing rows 7-39 down to 8-40) and populate it.
$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = 44545
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 100112032
$ws.Range("G7").Value = "Zapallo italiano"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("N7").Value = "$/caja 60 unidades"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 125
$ws.Range("Q7").Value = 60
$ws.Range("R7").Value = "Hortaliza"
